# TestDataMappingSheet.xlsx edit:
#   Insert a new "PlacementStability" mapping row into the
#   TestDataMappingSheet_SD sheet (right after the existing "Removal" row,
#   i.e. as the new row 123), pushing every row below it down by one.
#   Then fix up the AutoFilter range / _FilterDatabase defined name that
#   cover the table (they grow by one row), and leave the selection where
#   the user ended up after the edit (D126).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestDataMappingSheet_SD")

# --- Insert the new row, shifting rows 123:138 down to 124:139 -------------
$ws.Rows("123:123").Insert()

# --- Populate the newly inserted row ---------------------------------------
$ws.Range("A123").Value = "PlacementStability"
$ws.Range("B123").Value = "cares\Placement.xlsx"
$ws.Range("C123").Value = "PlacementStability"
$ws.Range("D123").Value = 1

# --- Grow the AutoFilter range from A1:E136 to A1:E137 ----------------------
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:E137").AutoFilter()

# --- Keep the _xlnm._FilterDatabase defined name in sync with the filter ---
$filterName = $wb.Names.Item("TestDataMappingSheet_SD!_FilterDatabase")
$filterName.RefersTo = "=TestDataMappingSheet_SD!`$A`$1:`$E`$137"

# --- Restore the active selection left after the edit -----------------------
$null = $ws.Range("D126").Select()
